$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 5 (extra student rows), keeping header + one data row
$ws.Range("A3:B5").EntireRow.Delete() | Out-Null

# Update remaining data row with new values
$ws.Range("A2").Value = 1815371
$ws.Range("B2").Value = "sol student 1"

# Update selection to match the new last-used cell
$ws.Range("B2").Select() | Out-Null
